$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing A14 timestamp with a more precise value
$ws.Range("A14").Value = 45863.87523748843

# Append new row 15 with the latest sensor reading
$ws.Range("A15").Value = 45863.9169284381
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 13.33
$ws.Range("E15").Value = 89.84999999999999
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 3.12
$ws.Range("H15").Value = "ENE"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "22:00:22"

# Match the date/time style used by column A in the existing rows
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat
